$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stock_list")

# New stock entries to insert, keyed by the worksheet row (in the ORIGINAL,
# pre-insert layout) they must be inserted immediately above. Processing
# these from the bottom of the sheet upwards means each insertion point is
# still valid (untouched by subsequent shifts) when we reach it.
# Each block is an ordered list of row hashtables (avoids PowerShell's
# array-unrolling surprises that a plain array-of-arrays runs into).

$blocks = @(
    @{ Before = 53; Rows = @(
            @{ A = "VFMO"; B = "Vanguard US Momentum Factor ETF"; C = "Exchange-traded fund (ETF)" },
            @{ A = "VZ";   B = "Verizon Communications Inc.";     C = "Telecommunications" }
        )
    },
    @{ Before = 45; Rows = @(
            @{ A = "SPGI"; B = "S&P Global";          C = "Financial" },
            @{ A = "SDGR"; B = "Schrödinger, Inc.";   C = "Chemical & Functional Materials" }
        )
    },
    @{ Before = 40; Rows = @(
            @{ A = "NVDA"; B = "NVIDIA"; C = "Semiconductors" }
        )
    },
    @{ Before = 30; Rows = @(
            @{ A = "KVUE"; B = "Kenvue"; C = "Consumer & Medical Goods" }
        )
    },
    @{ Before = 27; Rows = @(
            @{ A = "IVV";  B = "iShares Core S&P 500 ETF";  C = "Exchange-traded fund (ETF)" },
            @{ A = "MCHI"; B = "iShares MSCI China ETF";    C = "Exchange-traded fund (ETF)" },
            @{ A = "INDA"; B = "iShares MSCI India ETF";    C = "Exchange-traded fund (ETF)" },
            @{ A = "URTH"; B = "iShares MSCI World ETF";    C = "Exchange-traded fund (ETF)" }
        )
    },
    @{ Before = 25; Rows = @(
            @{ A = "GEV"; B = "GE Vernova Inc";        C = "Oil, Gas, & Energy" },
            @{ A = "HSY"; B = "The Hershey Company";   C = "Food & Fragrances" }
        )
    }
)

foreach ($block in $blocks) {
    $startRow = $block.Before
    $count = $block.Rows.Count

    # Insert `count` blank rows above $startRow.
    $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($startRow + $count - 1, 1)).EntireRow.Insert()

    # Fill the newly inserted blank rows with data.
    $r = $startRow
    foreach ($rowData in $block.Rows) {
        $ws.Cells.Item($r, 1).Value = $rowData.A
        $ws.Cells.Item($r, 2).Value = $rowData.B
        $ws.Cells.Item($r, 3).Value = $rowData.C
        $r = $r + 1
    }
}

$ws.Range("A27").Select()
